# Weekly data refresh: a new week's price record is inserted at row 44
# (pushing the existing rows 44-81 down to 45-82), matching the upstream
# "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 44, shifting rows 44:81 down
# to 45:82 (dimension grows from A1:R81 to A1:R82).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new week's record.
$ws.Range("A44").Value = 5
$ws.Range("B44").Value = 'Macroferia Regional de Talca'
$ws.Range("C44").Value = 'Maule'
$ws.Range("D44").Value = 44880
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = 100112040
$ws.Range("G44").Value = 'Cilantro'
$ws.Range("H44").Value = 'Sin especificar'
$ws.Range("I44").Value = 'Primera'
$ws.Range("J44").Value = 150
$ws.Range("K44").Value = 7000
$ws.Range("L44").Value = 7000
$ws.Range("M44").Value = 7000
$ws.Range("N44").Value = '$/caja 36 atados'
$ws.Range("O44").Value = 'Región del Maule'
$ws.Range("P44").Value = 194
$ws.Range("Q44").Value = 36
$ws.Range("R44").Value = 'Hortaliza'
